$wb = $excel.ActiveWorkbook

# Sheet ALC, row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1497.1666
$ws.Range("J32").Value = 1497.1666
$ws.Range("L32").Value = 1497.1666
$ws.Range("N32").Value = -2149.1666

# Sheet ALC, row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2496.6667
$ws.Range("J51").Value = 3500
$ws.Range("L51").Value = 3500
$ws.Range("N51").Value = -4468

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1555.4546
$ws.Range("I62").Value = 1610.4
$ws.Range("J62").Value = 1006
$ws.Range("K62").Value = 1610.4
$ws.Range("L62").Value = 1006
$ws.Range("M62").Value = -986.4000000000001
$ws.Range("N62").Value = -2254

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 1555.4546
$ws.Range("I65").Value = 1610.4
$ws.Range("J65").Value = 1006
$ws.Range("K65").Value = 8052
$ws.Range("L65").Value = 5030
$ws.Range("M65").Value = -4932
$ws.Range("N65").Value = -11270

# Sheet ALC, row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2290.7693
$ws.Range("I98").Value = 2472.5
$ws.Range("J98").Value = 2000
$ws.Range("K98").Value = 2472.5
$ws.Range("L98").Value = 2000
$ws.Range("M98").Value = -974.5
$ws.Range("N98").Value = -4996

# Sheet ALC, row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 11905708
$ws.Range("I112").Value = 837.1429000000001
$ws.Range("J112").Value = 14286682
$ws.Range("K112").Value = 2511.4287
$ws.Range("L112").Value = 42860046
$ws.Range("M112").Value = -1403.4287
$ws.Range("N112").Value = -42862262

# Sheet ALC, row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6730.952
$ws.Range("I116").Value = 12854.667
$ws.Range("J116").Value = 2138.1667
$ws.Range("K116").Value = 12854.667
$ws.Range("L116").Value = 2138.1667
$ws.Range("M116").Value = -9412.666999999999
$ws.Range("N116").Value = -9022.1667

# Sheet ALC, row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2290.7693
$ws.Range("I122").Value = 2472.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 7417.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -4967.5
$ws.Range("N122").Value = -10900

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1680.9667
$ws.Range("J137").Value = 1648.8
$ws.Range("L137").Value = 4946.4
$ws.Range("N137").Value = -10046.4

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5299.8833
$ws.Range("I32").Value = 4329.362
$ws.Range("K32").Value = 4329.362
$ws.Range("M32").Value = -4042.362

# Sheet ARM, row 37
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 27666.666
$ws.Range("I37").Value = 9000
$ws.Range("J37").Value = 37000
$ws.Range("K37").Value = 9000
$ws.Range("L37").Value = 37000
$ws.Range("M37").Value = -8727
$ws.Range("N37").Value = -37546

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6274.609
$ws.Range("I74").Value = 8844
$ws.Range("J74").Value = 2277.7778
$ws.Range("K74").Value = 8844
$ws.Range("L74").Value = 2277.7778
$ws.Range("M74").Value = -7970
$ws.Range("N74").Value = -4025.7778

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6274.609
$ws.Range("I77").Value = 8844
$ws.Range("J77").Value = 2277.7778
$ws.Range("K77").Value = 44220
$ws.Range("L77").Value = 11388.889
$ws.Range("M77").Value = -39852
$ws.Range("N77").Value = -20124.889

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6653.6772
$ws.Range("I132").Value = 2034.7894
$ws.Range("J132").Value = 13966.917
$ws.Range("K132").Value = 6104.3682
$ws.Range("L132").Value = 41900.751
$ws.Range("M132").Value = -3574.3682
$ws.Range("N132").Value = -46960.751

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1433.6666
$ws.Range("I94").Value = 1297.8649
$ws.Range("J94").Value = 1890.4546
$ws.Range("K94").Value = 1297.8649
$ws.Range("L94").Value = 1890.4546
$ws.Range("M94").Value = -846.8649
$ws.Range("N94").Value = -2792.4546

# Sheet BSM, row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1300.5264
$ws.Range("I107").Value = 1258.375
$ws.Range("J107").Value = 1525.3334
$ws.Range("K107").Value = 1258.375
$ws.Range("L107").Value = 1525.3334
$ws.Range("M107").Value = 661.625
$ws.Range("N107").Value = -5365.3334

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8411.972
$ws.Range("I31").Value = 1698.25
$ws.Range("J31").Value = 11914.782
$ws.Range("K31").Value = 1698.25
$ws.Range("L31").Value = 11914.782
$ws.Range("M31").Value = -1403.25
$ws.Range("N31").Value = -12504.782

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8411.972
$ws.Range("I34").Value = 1698.25
$ws.Range("J34").Value = 11914.782
$ws.Range("K34").Value = 1698.25
$ws.Range("L34").Value = 11914.782
$ws.Range("M34").Value = -1496.25
$ws.Range("N34").Value = -12318.782

# Sheet CRP, row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 2593.32
$ws.Range("I94").Value = 1991
$ws.Range("J94").Value = 3149.3076
$ws.Range("K94").Value = 1991
$ws.Range("L94").Value = 3149.3076
$ws.Range("M94").Value = -1540
$ws.Range("N94").Value = -4051.3076

# Sheet CRP, row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 796
$ws.Range("I107").Value = 696.4666999999999
$ws.Range("J107").Value = 902.6429000000001
$ws.Range("K107").Value = 696.4666999999999
$ws.Range("L107").Value = 902.6429000000001
$ws.Range("M107").Value = 1223.5333
$ws.Range("N107").Value = -4742.6429

# Sheet GSM, row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 66667892
$ws.Range("I113").Value = 100001080
$ws.Range("J113").Value = 1520
$ws.Range("K113").Value = 100001080
$ws.Range("L113").Value = 1520
$ws.Range("M113").Value = -99998910
$ws.Range("N113").Value = -5860

# Sheet GSM, row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6857.6
$ws.Range("I126").Value = 9441.308000000001
$ws.Range("K126").Value = 28323.924
$ws.Range("M126").Value = -25853.924

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2929.0962
$ws.Range("I132").Value = 2792.276
$ws.Range("J132").Value = 3101.6086
$ws.Range("K132").Value = 8376.828
$ws.Range("L132").Value = 9304.825800000001
$ws.Range("M132").Value = -5846.828
$ws.Range("N132").Value = -14364.8258

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 22223462
$ws.Range("I46").Value = 41667680
$ws.Range("J46").Value = 1499.7142
$ws.Range("K46").Value = 41667680
$ws.Range("L46").Value = 1499.7142
$ws.Range("M46").Value = -41667492
$ws.Range("N46").Value = -1875.7142

# Sheet LTW, row 92
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 22500
$ws.Range("J92").Value = 22500
$ws.Range("L92").Value = 22500
$ws.Range("N92").Value = -27492

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11408691
$ws.Range("I132").Value = 18846882
$ws.Range("J132").Value = 3465.0667
$ws.Range("K132").Value = 56540646
$ws.Range("L132").Value = 10395.2001
$ws.Range("M132").Value = -56538116
$ws.Range("N132").Value = -15455.2001

# Sheet WVR, row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3170.6155
$ws.Range("I122").Value = 1534
$ws.Range("J122").Value = 4573.4287
$ws.Range("K122").Value = 4602
$ws.Range("L122").Value = 13720.2861
$ws.Range("M122").Value = -2152
$ws.Range("N122").Value = -18620.2861
